$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row: Cost1 (F2) for row 2
$ws.Range("F2").Value = 70

# Row 5: Calizon Dike
$ws.Range("A5").Value = $true
$ws.Range("B5").Value = "Calizon Dike"
$ws.Range("C5").Value = 14.9136800407707
$ws.Range("D5").Value = 120.755871075221
$ws.Range("E5").Value = 126
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 252
$ws.Range("H5").Value = 7560000
$ws.Range("I5").Value = $false
$ws.Range("J5").Value = $false
$ws.Range("K5").Value = $false
$ws.Range("L5").Value = "Built"

# Row 6: Frances E.C.
$ws.Range("A6").Value = $true
$ws.Range("B6").Value = "Frances E.C."
$ws.Range("C6").Value = 14.9194611702998
$ws.Range("D6").Value = 120.762172685224
$ws.Range("E6").Value = 150
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 300
$ws.Range("H6").Value = 9000000
$ws.Range("I6").Value = $true
$ws.Range("J6").Value = $true
$ws.Range("K6").Value = $true
$ws.Range("L6").Value = "Built"
